$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so numeric-looking values are not
# auto-converted to numbers by Excel.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '69.911.16'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '3.535.73'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '602.86'
$ws.Range('E5').Value = '  -2.15%  '
$ws.Range('D6').Value = '196.70'
$ws.Range('E6').Value = '  +5.59%  '
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -2.93%  '
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '54.02'
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '0.0000302'
$ws.Range('E12').Value = '  -2.71%  '
$ws.Range('D13').Value = '9.54'
$ws.Range('E13').Value = '  -0.13%  '
$ws.Range('D14').Value = '4.084.44'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').Value = '599.05'
$ws.Range('E15').Value = '  -4.64%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '70.100.67'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '19.15'
$ws.Range('E17').Value = '  +0.83%  '
$ws.Range('D18').Value = '12.66'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('D19').Value = '3.525.17'
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').Value = '18.16'
$ws.Range('E22').Value = '  +3.07%  '
$ws.Range('D23').Value = '5.30'
$ws.Range('E23').Value = '  +6.90%  '
$ws.Range('D24').Value = '103.37'
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').Value = '4.61'
$ws.Range('E25').Value = '  -2.57%  '
$ws.Range('E26').Value = '  +2.25%  '
$ws.Range('D27').Value = '10.94'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '9.68'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').Value = '33.57'
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('D30').Value = '4.50'
$ws.Range('E30').Value = '  +20.30%  '
$ws.Range('D31').Value = '7.10'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  +3.43%  '
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('D34').Value = '63.50'
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('D35').Value = '0.0₃0830'
$ws.Range('E35').Value = '  +5.78%  '
$ws.Range('D36').Value = '3.741.94'
$ws.Range('E36').Value = '  +4.23%  '
$ws.Range('E37').Value = '  -5.72%  '
$ws.Range('E38').Value = '  -0.20%  '
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('D40').Value = '3.60'
$ws.Range('E40').Value = '  +1.57%  '
$ws.Range('D41').Value = '36.81'
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('D42').Value = '495.14'
$ws.Range('E42').Value = '  -7.04%  '
$ws.Range('E43').Value = '  -1.05%  '
$ws.Range('D44').Value = '0.0456'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('E45').Value = '  -2.37%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '3.33'
$ws.Range('E46').Value = '  -1.90%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').Value = '2.83'
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('D49').Value = '8.67'
$ws.Range('E49').Value = '  -5.89%  '
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('E51').Value = '  +12.20%  '
